$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns keep their text representation
# (values like "1.00", "0.999", "6.89" must not be coerced to numbers)
$ws.Range("D2:E51").NumberFormat = "@"

$updates = @{
    "D2" = "57.460.80"
    "E2" = "  -4.62%  "
    "D3" = "3.116.52"
    "E3" = "  -5.86%  "
    "D4" = "1.00"
    "E4" = "  +0.08%  "
    "D5" = "519.82"
    "E5" = "  -7.05%  "
    "D6" = "133.95"
    "E6" = "  -5.67%  "
    "E7" = "  -0.11%  "
    "D8" = "3.117.80"
    "E8" = "  -5.78%  "
    "D9" = "0.443"
    "E9" = "  -6.55%  "
    "E10" = "  -8.45%  "
    "E11" = "  -8.43%  "
    "D12" = "0.381"
    "E12" = "  -6.27%  "
    "D13" = "3.656.75"
    "E13" = "  -5.74%  "
    "E14" = "  -2.32%  "
    "D15" = "25.24"
    "E15" = "  -5.92%  "
    "D16" = "3.118.51"
    "E16" = "  -5.66%  "
    "D17" = "57.469.33"
    "E17" = "  -4.62%  "
    "D18" = "0.0000150"
    "E18" = "  -9.29%  "
    "D19" = "5.73"
    "E19" = "  -7.14%  "
    "D20" = "12.90"
    "E20" = "  -10.59%  "
    "D21" = "7.92"
    "E21" = "  -8.26%  "
    "D22" = "340.71"
    "E22" = "  -9.16%  "
    "E23" = "  -0.08%  "
    "D24" = "68.14"
    "E24" = "  -8.16%  "
    "D25" = "0.501"
    "E25" = "  -7.59%  "
    "D26" = "3.244.73"
    "E26" = "  -5.83%  "
    "E27" = "  -3.83%  "
    "D28" = "0.999"
    "E28" = "  -0.19%  "
    "E29" = "  -9.53%  "
    "E30" = "  -0.27%  "
    "D31" = "6.74"
    "E31" = "  -6.86%  "
    "B32" = "PancakeSwap"
    "C32" = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
    "D32" = "1.86"
    "E32" = "  -8.51%  "
    "B33" = "InternetComputer(DFINITY)"
    "C33" = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
    "D33" = "6.89"
    "E33" = "  -9.87%  "
    "E34" = "  -3.46%  "
    "D35" = "21.37"
    "E35" = "  -5.12%  "
    "D36" = "158.42"
    "E36" = "  -4.47%  "
    "D37" = "4.74"
    "E37" = "  -7.84%  "
    "D38" = "6.12"
    "E38" = "  -9.18%  "
    "D39" = "1.36"
    "E39" = "  -11.14%  "
    "D40" = "25.11"
    "E40" = "  -6.17%  "
    "D41" = "0.0684"
    "E41" = "  -7.21%  "
    "D42" = "3.149.55"
    "E42" = "  -5.69%  "
    "D43" = "40.27"
    "E43" = "  -4.07%  "
    "D44" = "0.680"
    "E44" = "  -9.65%  "
    "B45" = "ONDO"
    "C45" = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
    "D45" = "1.06"
    "E45" = "  -4.63%  "
    "B46" = "Filecoin"
    "C46" = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
    "D46" = "3.88"
    "E46" = "  -7.45%  "
    "D47" = "1.00"
    "E47" = "  +0.15%  "
    "D49" = "2.247.24"
    "E49" = "  -4.79%  "
    "D50" = "6.14"
    "E50" = "  -5.87%  "
    "D51" = "19.80"
    "E51" = "  -7.04%  "
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
